$wb = $excel.ActiveWorkbook
$ranking = $wb.Worksheets.Item("Ranking")
$passFail = $wb.Worksheets.Item("Pass Fail")

# Data edit: Full Tests project "schedule" - Overlap/Other scores changed
$ranking.Range("D14").Value = 2
$ranking.Range("E14").Value = 3

# Rename "Pass Rate" label to "Average Pass Rate" on the Pass Fail sheet
$passFail.Range("A17").Value = "Average Pass Rate"
$passFail.Range("A31").Value = "Average Pass Rate"
